$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A152").Value = "IMX-USD"
$ws.Range("A153").Value = "GRT-USD"
